$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the floating point timestamp on the last existing row (A35).
$ws.Range("A35").Value = 44348.86069004977

# Append the new data row (row 36).
$ws.Range("A36").Value = 44349.85674851396
$ws.Range("B36").Value = 74934
$ws.Range("C36").Value = 63172
$ws.Range("D36").Value = 3321
$ws.Range("E36").Value = 2082
$ws.Range("F36").Value = 1468
$ws.Range("G36").Value = 19618
$ws.Range("H36").Value = 1394
$ws.Range("I36").Value = 871
$ws.Range("J36").Value = 206

# Match the date-formatted style used by the rest of column A.
$ws.Range("A36").NumberFormat = $ws.Range("A35").NumberFormat
